$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1369.6086
$ws.Range("I15").Value = 1369.6086
$ws.Range("K15").Value = 4108.825800000001
$ws.Range("M15").Value = -3939.825800000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 3519.7856
$ws.Range("I98").Value = 3542.2856
$ws.Range("J98").Value = 3497.2856
$ws.Range("K98").Value = 3542.2856
$ws.Range("L98").Value = 3497.2856
$ws.Range("M98").Value = -2044.2856
$ws.Range("N98").Value = -6493.2856

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 617.28125
$ws.Range("J107").Value = 574.6667
$ws.Range("L107").Value = 574.6667
$ws.Range("N107").Value = -4414.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2049
$ws.Range("I112").Value = 1500
$ws.Range("J112").Value = 2186.25
$ws.Range("K112").Value = 4500
$ws.Range("L112").Value = 6558.75
$ws.Range("M112").Value = -3392
$ws.Range("N112").Value = -8774.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 3519.7856
$ws.Range("I122").Value = 3542.2856
$ws.Range("J122").Value = 3497.2856
$ws.Range("K122").Value = 10626.8568
$ws.Range("L122").Value = 10491.8568
$ws.Range("M122").Value = -8176.856800000001
$ws.Range("N122").Value = -15391.8568

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2871.8809
$ws.Range("I132").Value = 2769.4167
$ws.Range("K132").Value = 8308.250100000001
$ws.Range("M132").Value = -5778.250100000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1535
$ws.Range("I135").Value = 437
$ws.Range("K135").Value = 3933
$ws.Range("M135").Value = -1398

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2234.2144
$ws.Range("I137").Value = 1770.909
$ws.Range("J137").Value = 3933
$ws.Range("K137").Value = 5312.727000000001
$ws.Range("L137").Value = 11799
$ws.Range("M137").Value = -2762.727000000001
$ws.Range("N137").Value = -16899

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3211.7646
$ws.Range("I138").Value = 1499
$ws.Range("J138").Value = 3440.1333
$ws.Range("K138").Value = 4497
$ws.Range("L138").Value = 10320.3999
$ws.Range("M138").Value = 643
$ws.Range("N138").Value = -20600.3999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2689.3333
$ws.Range("I45").Value = 2458.4285
$ws.Range("K45").Value = 2458.4285
$ws.Range("M45").Value = -2081.4285

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1380.9231
$ws.Range("I74").Value = 1262.6666
$ws.Range("J74").Value = 2800
$ws.Range("K74").Value = 1262.6666
$ws.Range("L74").Value = 2800
$ws.Range("M74").Value = -388.6666
$ws.Range("N74").Value = -4548

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1380.9231
$ws.Range("I77").Value = 1262.6666
$ws.Range("J77").Value = 2800
$ws.Range("K77").Value = 6313.333000000001
$ws.Range("L77").Value = 14000
$ws.Range("M77").Value = -1945.333000000001
$ws.Range("N77").Value = -22736

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2069.889
$ws.Range("I102").Value = 2069.889
$ws.Range("K102").Value = 2069.889
$ws.Range("M102").Value = -447.8890000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 4442
$ws.Range("I110").Value = 2768.3333
$ws.Range("K110").Value = 2768.3333
$ws.Range("M110").Value = -723.3332999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3215.9678
$ws.Range("I132").Value = 3008.4482
$ws.Range("J132").Value = 6225
$ws.Range("K132").Value = 9025.3446
$ws.Range("L132").Value = 18675
$ws.Range("M132").Value = -6495.3446
$ws.Range("N132").Value = -23735

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2064.3684
$ws.Range("I86").Value = 2144.8
$ws.Range("J86").Value = 1762.75
$ws.Range("K86").Value = 2144.8
$ws.Range("L86").Value = 1762.75
$ws.Range("M86").Value = -1021.8
$ws.Range("N86").Value = -4008.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2064.3684
$ws.Range("I89").Value = 2144.8
$ws.Range("J89").Value = 1762.75
$ws.Range("K89").Value = 10724
$ws.Range("L89").Value = 8813.75
$ws.Range("M89").Value = -5108
$ws.Range("N89").Value = -20045.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3521.647
$ws.Range("J31").Value = 3281.25
$ws.Range("L31").Value = 3281.25
$ws.Range("N31").Value = -3871.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3521.647
$ws.Range("J34").Value = 3281.25
$ws.Range("L34").Value = 3281.25
$ws.Range("N34").Value = -3685.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 48333.332
$ws.Range("J64").Value = 48333.332
$ws.Range("L64").Value = 48333.332
$ws.Range("N64").Value = -48829.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H67").Value = 48333.332
$ws.Range("J67").Value = 48333.332
$ws.Range("L67").Value = 48333.332
$ws.Range("N67").Value = -50049.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H80").Value = 10799.8
$ws.Range("J80").Value = 10799.8
$ws.Range("L80").Value = 10799.8
$ws.Range("N80").Value = -13045.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H83").Value = 10799.8
$ws.Range("J83").Value = 10799.8
$ws.Range("L83").Value = 32399.4
$ws.Range("N83").Value = -43631.39999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H103").Value = 9474.666999999999
$ws.Range("I103").Value = 9474.666999999999
$ws.Range("K103").Value = 9474.666999999999
$ws.Range("M103").Value = -8302.666999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 474.1875
$ws.Range("I107").Value = 245.83333
$ws.Range("J107").Value = 611.2
$ws.Range("K107").Value = 245.83333
$ws.Range("L107").Value = 611.2
$ws.Range("M107").Value = 1674.16667
$ws.Range("N107").Value = -4451.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H109").Value = 57665.668
$ws.Range("J109").Value = 57665.668
$ws.Range("L109").Value = 57665.668
$ws.Range("N109").Value = -59745.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2848.1667
$ws.Range("I132").Value = 3257.8
$ws.Range("K132").Value = 9773.400000000001
$ws.Range("M132").Value = -7243.400000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2143.5454
$ws.Range("I134").Value = 2274
$ws.Range("J134").Value = 1700
$ws.Range("K134").Value = 6822
$ws.Range("L134").Value = 5100
$ws.Range("M134").Value = -4287
$ws.Range("N134").Value = -10170

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1378.5758
$ws.Range("I4").Value = 876.9545000000001
$ws.Range("J4").Value = 2381.818
$ws.Range("K4").Value = 2630.8635
$ws.Range("L4").Value = 7145.454000000001
$ws.Range("M4").Value = -2518.8635
$ws.Range("N4").Value = -7369.454000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 292.25
$ws.Range("I107").Value = 292.25
$ws.Range("K107").Value = 292.25
$ws.Range("M107").Value = 1627.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1175.625
$ws.Range("I132").Value = 986.4286
$ws.Range("K132").Value = 2959.2858
$ws.Range("M132").Value = -429.2857999999997

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2310
$ws.Range("I16").Value = 2636.2222
$ws.Range("J16").Value = 1331.3334
$ws.Range("K16").Value = 2636.2222
$ws.Range("L16").Value = 1331.3334
$ws.Range("M16").Value = -2466.2222
$ws.Range("N16").Value = -1671.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2731
$ws.Range("I46").Value = 1195
$ws.Range("J46").Value = 3499
$ws.Range("K46").Value = 1195
$ws.Range("L46").Value = 3499
$ws.Range("M46").Value = -1007
$ws.Range("N46").Value = -3875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 838.4167
$ws.Range("I55").Value = 840.1111
$ws.Range("K55").Value = 840.1111
$ws.Range("M55").Value = -667.1111

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3541.1428
$ws.Range("I136").Value = 3541.1428
$ws.Range("K136").Value = 10623.4284
$ws.Range("M136").Value = -8073.428400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 119999
$ws.Range("J109").Value = 119999
$ws.Range("L109").Value = 119999
$ws.Range("N109").Value = -122773

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 481.25
$ws.Range("I132").Value = 437.5
$ws.Range("J132").Value = 525
$ws.Range("K132").Value = 1312.5
$ws.Range("L132").Value = 1575
$ws.Range("N132").Value = -6635

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3972.7778
$ws.Range("I136").Value = 3972.7778
$ws.Range("K136").Value = 11918.3334
$ws.Range("M136").Value = -9368.3334

# Special cases: cell removals and additions
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M99").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("N98").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M132").Value = 1217.5
